$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 0.3201913333333333
$ws.Range("H2").Value = 0.960574
$ws.Range("I2").Value = 0.03939146858412543
$ws.Range("J2").Value = 0.03939146858412543
$ws.Range("M2").Value = 4.277890333333334
$ws.Range("N2").Value = 12.833671
$ws.Range("O2").Value = 0.04123357425337639
$ws.Range("P2").Value = 0.04123357425337638
$ws.Range("Q2").Value = 1.369743409683778
$ws.Range("R2").Value = 12.327690687154
$ws.Range("S2").Value = 0.001624251044813079
$ws.Range("T2").Value = 0.001624251044813079

# Row 3
$ws.Range("G3").Value = 0.3201913333333333
$ws.Range("H3").Value = 0.960574
$ws.Range("I3").Value = 0.03939146858412543
$ws.Range("J3").Value = 0.03939146858412543
$ws.Range("O3").Value = 0.4451428460610328
$ws.Range("P3").Value = 0.4451428460610327
$ws.Range("Q3").Value = 14.78725749102511
$ws.Range("R3").Value = 133.085317419226
$ws.Range("S3").Value = 0.01753483043606136
$ws.Range("T3").Value = 0.01753483043606135

# Row 4
$ws.Range("G4").Value = 0.3201913333333333
$ws.Range("H4").Value = 0.960574
$ws.Range("I4").Value = 0.03939146858412543
$ws.Range("J4").Value = 0.03939146858412543
$ws.Range("M4").Value = 8.558147333333332
$ws.Range("N4").Value = 25.674442
$ws.Range("O4").Value = 0.08248996024761777
$ws.Range("P4").Value = 0.08248996024761777
$ws.Range("Q4").Value = 2.740244605523111
$ws.Range("R4").Value = 24.662201449708
$ws.Range("S4").Value = 0.003249400677599791
$ws.Range("T4").Value = 0.003249400677599791

# Row 5
$ws.Range("G5").Value = 0.3201913333333333
$ws.Range("H5").Value = 0.960574
$ws.Range("I5").Value = 0.03939146858412543
$ws.Range("J5").Value = 0.03939146858412543
$ws.Range("M5").Value = 44.72914066666667
$ws.Range("N5").Value = 134.187422
$ws.Range("O5").Value = 0.4311336194379731
$ws.Range("P5").Value = 0.431133619437973
$ws.Range("Q5").Value = 14.32188318891422
$ws.Range("R5").Value = 128.896948700228
$ws.Range("S5").Value = 0.01698298642565121
$ws.Range("T5").Value = 0.0169829864256512

# Row 6
$ws.Range("I6").Value = 0.2346323697636092
$ws.Range("J6").Value = 0.2346323697636091
$ws.Range("M6").Value = 4.277890333333334
$ws.Range("N6").Value = 12.833671
$ws.Range("O6").Value = 0.04123357425337639
$ws.Range("P6").Value = 0.04123357425337638
$ws.Range("Q6").Value = 8.158775332172002
$ws.Range("R6").Value = 73.42897798954802
$ws.Range("S6").Value = 0.009674731240893443
$ws.Range("T6").Value = 0.00967473124089344

# Row 7
$ws.Range("I7").Value = 0.2346323697636092
$ws.Range("J7").Value = 0.2346323697636091
$ws.Range("O7").Value = 0.4451428460610328
$ws.Range("P7").Value = 0.4451428460610327
$ws.Range("S7").Value = 0.1044449208546176
$ws.Range("T7").Value = 0.1044449208546176

# Row 8
$ws.Range("I8").Value = 0.2346323697636092
$ws.Range("J8").Value = 0.2346323697636091
$ws.Range("M8").Value = 8.558147333333332
$ws.Range("N8").Value = 25.674442
$ws.Range("O8").Value = 0.08248996024761777
$ws.Range("P8").Value = 0.08248996024761777
$ws.Range("Q8").Value = 16.322064361544
$ws.Range("R8").Value = 146.898579253896
$ws.Range("S8").Value = 0.01935481485460448
$ws.Range("T8").Value = 0.01935481485460447

# Row 9
$ws.Range("I9").Value = 0.2346323697636092
$ws.Range("J9").Value = 0.2346323697636091
$ws.Range("M9").Value = 44.72914066666667
$ws.Range("N9").Value = 134.187422
$ws.Range("O9").Value = 0.4311336194379731
$ws.Range("P9").Value = 0.431133619437973
$ws.Range("Q9").Value = 85.307238162904
$ws.Range("R9").Value = 767.765143466136
$ws.Range("S9").Value = 0.1011579028134937
$ws.Range("T9").Value = 0.1011579028134936

# Row 10
$ws.Range("G10").Value = 5.780535
$ws.Range("H10").Value = 17.341605
$ws.Range("I10").Value = 0.7111490510422025
$ws.Range("J10").Value = 0.7111490510422023
$ws.Range("M10").Value = 4.277890333333334
$ws.Range("N10").Value = 12.833671
$ws.Range("O10").Value = 0.04123357425337639
$ws.Range("P10").Value = 0.04123357425337638
$ws.Range("Q10").Value = 24.72849479799501
$ws.Range("R10").Value = 222.556453181955
$ws.Range("S10").Value = 0.02932321720136681
$ws.Range("T10").Value = 0.0293232172013668

# Row 11
$ws.Range("G11").Value = 5.780535
$ws.Range("H11").Value = 17.341605
$ws.Range("I11").Value = 0.7111490510422025
$ws.Range("J11").Value = 0.7111490510422023
$ws.Range("O11").Value = 0.4451428460610328
$ws.Range("P11").Value = 0.4451428460610327
$ws.Range("Q11").Value = 266.959941079655
$ws.Range("R11").Value = 2402.639469716895
$ws.Range("S11").Value = 0.3165629125545287
$ws.Range("T11").Value = 0.3165629125545286

# Row 12
$ws.Range("G12").Value = 5.780535
$ws.Range("H12").Value = 17.341605
$ws.Range("I12").Value = 0.7111490510422025
$ws.Range("J12").Value = 0.7111490510422023
$ws.Range("M12").Value = 8.558147333333332
$ws.Range("N12").Value = 25.674442
$ws.Range("O12").Value = 0.08248996024761777
$ws.Range("P12").Value = 0.08248996024761777
$ws.Range("Q12").Value = 49.47067019549
$ws.Range("R12").Value = 445.23603175941
$ws.Range("S12").Value = 0.05866265695060239
$ws.Range("T12").Value = 0.05866265695060237

# Row 13
$ws.Range("G13").Value = 5.780535
$ws.Range("H13").Value = 17.341605
$ws.Range("I13").Value = 0.7111490510422025
$ws.Range("J13").Value = 0.7111490510422023
$ws.Range("M13").Value = 44.72914066666667
$ws.Range("N13").Value = 134.187422
$ws.Range("O13").Value = 0.4311336194379731
$ws.Range("P13").Value = 0.431133619437973
$ws.Range("Q13").Value = 258.55836314359
$ws.Range("R13").Value = 2327.02526829231
$ws.Range("S13").Value = 0.3066002643357046
$ws.Range("T13").Value = 0.3066002643357045

# Row 14
$ws.Range("E14").Value = 2
$ws.Range("F14").Value = 0.6666666666666666
$ws.Range("G14").Value = 0.1205213333333333
$ws.Range("H14").Value = 0.361564
$ws.Range("I14").Value = 0.01482711061006308
$ws.Range("J14").Value = 0.01482711061006307
$ws.Range("M14").Value = 4.277890333333334
$ws.Range("N14").Value = 12.833671
$ws.Range("O14").Value = 0.04123357425337639
$ws.Range("P14").Value = 0.04123357425337638
$ws.Range("Q14").Value = 0.5155770468271113
$ws.Range("R14").Value = 4.640193421444001
$ws.Range("S14").Value = 0.0006113747663030608
$ws.Range("T14").Value = 0.0006113747663030605

# Row 15
$ws.Range("E15").Value = 2
$ws.Range("F15").Value = 0.6666666666666666
$ws.Range("G15").Value = 0.1205213333333333
$ws.Range("H15").Value = 0.361564
$ws.Range("I15").Value = 0.01482711061006308
$ws.Range("J15").Value = 0.01482711061006307
$ws.Range("O15").Value = 0.4451428460610328
$ws.Range("P15").Value = 0.4451428460610327
$ws.Range("Q15").Value = 5.565984471248444
$ws.Range("R15").Value = 50.09386024123599
$ws.Range("S15").Value = 0.006600182215825213
$ws.Range("T15").Value = 0.006600182215825211

# Row 16
$ws.Range("E16").Value = 2
$ws.Range("F16").Value = 0.6666666666666666
$ws.Range("G16").Value = 0.1205213333333333
$ws.Range("H16").Value = 0.361564
$ws.Range("I16").Value = 0.01482711061006308
$ws.Range("J16").Value = 0.01482711061006307
$ws.Range("M16").Value = 8.558147333333332
$ws.Range("N16").Value = 25.674442
$ws.Range("O16").Value = 0.08248996024761777
$ws.Range("P16").Value = 0.08248996024761777
$ws.Range("Q16").Value = 1.031439327476444
$ws.Range("R16").Value = 9.282953947288
$ws.Range("S16").Value = 0.001223087764811135
$ws.Range("T16").Value = 0.001223087764811135

# Row 17
$ws.Range("E17").Value = 2
$ws.Range("F17").Value = 0.6666666666666666
$ws.Range("G17").Value = 0.1205213333333333
$ws.Range("H17").Value = 0.361564
$ws.Range("I17").Value = 0.01482711061006308
$ws.Range("J17").Value = 0.01482711061006307
$ws.Range("M17").Value = 44.72914066666667
$ws.Range("N17").Value = 134.187422
$ws.Range("O17").Value = 0.4311336194379731
$ws.Range("P17").Value = 0.431133619437973
$ws.Range("Q17").Value = 5.390815672000889
$ws.Range("R17").Value = 48.517341048008
$ws.Range("S17").Value = 0.006392465863123667
$ws.Range("T17").Value = 0.006392465863123665
